# Apply resume edits per commit "changes to general resume"

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Title: FRONT-END DEVELOPER -> FULL-STACK WEB DEVELOPER
Replace-Text "FRONT-END DEVELOPER" "FULL-STACK WEB DEVELOPER"

# 2. Summary paragraph
Replace-Text "Detail-oriented Frontend Developer with a strong background in UI/UX design principles and front-end technologies. Excited to collaborate with cross-functional teams to deliver high-quality user interfaces and user experiences." "Innovative Full Stack Web Developer with experience in React and Node.js, specializing in building scalable web applications. Proven track record of delivering end-to-end solutions with a focus on user experience and performance. Eager to expand expertise into smart contract development and DevOps practices."

# 3. Version control: Git/Github -> Git
Replace-Text ": Git/Github" ": Git"

# 4. Additional Skills bullet list
Replace-Text "Detail-oriented" "Quick learner with a keen interest in new technologies"
Replace-Text "Time efficient" "Strong problem-solving skills and attention to detail"
Replace-Text "Flexibility and adaptability" "Excellent communication and interpersonal skills"
Replace-Text "Strong communication and interpersonal skills" "Highly flexible and adaptable"
